$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename SundayAM header (B1) to "Not Full" - a not-enough-staff warning column
$ws.Cells.Item(1, 2).Value = 'Not Full'

# Re-populate the shift grid (B2:O32) with the reshuffled schedule
$ws.Cells.Item(2, 2).Value = '-'
$ws.Cells.Item(2, 3).Value = '4:00: Bar'
$ws.Cells.Item(2, 4).Value = '-'
$ws.Cells.Item(2, 5).Value = '-'
$ws.Cells.Item(2, 6).Value = '-'
$ws.Cells.Item(2, 7).Value = '-'
$ws.Cells.Item(2, 8).Value = '-'
$ws.Cells.Item(2, 9).Value = '-'
$ws.Cells.Item(2, 10).Value = '-'
$ws.Cells.Item(2, 11).Value = '4:00: Bar'
$ws.Cells.Item(2, 12).Value = '-'
$ws.Cells.Item(2, 13).Value = '4:00: Bar'
$ws.Cells.Item(2, 14).Value = '-'
$ws.Cells.Item(2, 15).Value = '4:00: Bar'
$ws.Cells.Item(3, 2).Value = '-'
$ws.Cells.Item(3, 3).Value = '-'
$ws.Cells.Item(3, 4).Value = '10:00: Bar'
$ws.Cells.Item(3, 5).Value = '-'
$ws.Cells.Item(3, 6).Value = '10:00: Bar'
$ws.Cells.Item(3, 7).Value = '-'
$ws.Cells.Item(3, 8).Value = '-'
$ws.Cells.Item(3, 9).Value = '4:00: Bar'
$ws.Cells.Item(3, 10).Value = '10:00: Bar'
$ws.Cells.Item(3, 11).Value = '-'
$ws.Cells.Item(3, 12).Value = '10:00: Bar'
$ws.Cells.Item(3, 13).Value = '-'
$ws.Cells.Item(3, 14).Value = '-'
$ws.Cells.Item(3, 15).Value = '-'
$ws.Cells.Item(4, 2).Value = '-'
$ws.Cells.Item(4, 3).Value = '-'
$ws.Cells.Item(4, 4).Value = '-'
$ws.Cells.Item(4, 5).Value = '-'
$ws.Cells.Item(4, 6).Value = '-'
$ws.Cells.Item(4, 7).Value = '-'
$ws.Cells.Item(4, 8).Value = '-'
$ws.Cells.Item(4, 9).Value = '4:00: Bar'
$ws.Cells.Item(4, 10).Value = '-'
$ws.Cells.Item(4, 11).Value = '4:00: Bar'
$ws.Cells.Item(4, 12).Value = '-'
$ws.Cells.Item(4, 13).Value = '4:00: Bar'
$ws.Cells.Item(4, 14).Value = '10:00: Bar'
$ws.Cells.Item(4, 15).Value = '4:00: Bar'
$ws.Cells.Item(5, 2).Value = '-'
$ws.Cells.Item(5, 3).Value = '4:00: Bar'
$ws.Cells.Item(5, 4).Value = '-'
$ws.Cells.Item(5, 5).Value = '4:00: Bar'
$ws.Cells.Item(5, 6).Value = '-'
$ws.Cells.Item(5, 7).Value = '4:00: Bar'
$ws.Cells.Item(5, 8).Value = '10:00: Bar'
$ws.Cells.Item(5, 9).Value = '-'
$ws.Cells.Item(5, 10).Value = '-'
$ws.Cells.Item(5, 11).Value = '-'
$ws.Cells.Item(5, 12).Value = '-'
$ws.Cells.Item(5, 13).Value = '-'
$ws.Cells.Item(5, 14).Value = '-'
$ws.Cells.Item(5, 15).Value = '-'
$ws.Cells.Item(6, 2).Value = '-'
$ws.Cells.Item(6, 3).Value = '-'
$ws.Cells.Item(6, 4).Value = '-'
$ws.Cells.Item(6, 5).Value = '-'
$ws.Cells.Item(6, 6).Value = '-'
$ws.Cells.Item(6, 7).Value = '-'
$ws.Cells.Item(6, 8).Value = '-'
$ws.Cells.Item(6, 9).Value = '-'
$ws.Cells.Item(6, 10).Value = '-'
$ws.Cells.Item(6, 11).Value = '-'
$ws.Cells.Item(6, 12).Value = '-'
$ws.Cells.Item(6, 13).Value = '-'
$ws.Cells.Item(6, 14).Value = '-'
$ws.Cells.Item(6, 15).Value = '-'
$ws.Cells.Item(7, 2).Value = '-'
$ws.Cells.Item(7, 3).Value = '-'
$ws.Cells.Item(7, 4).Value = '-'
$ws.Cells.Item(7, 5).Value = '4:00: Bar'
$ws.Cells.Item(7, 6).Value = '-'
$ws.Cells.Item(7, 7).Value = '4:00: Bar'
$ws.Cells.Item(7, 8).Value = '-'
$ws.Cells.Item(7, 9).Value = '-'
$ws.Cells.Item(7, 10).Value = '-'
$ws.Cells.Item(7, 11).Value = '-'
$ws.Cells.Item(7, 12).Value = '-'
$ws.Cells.Item(7, 13).Value = '-'
$ws.Cells.Item(7, 14).Value = '-'
$ws.Cells.Item(7, 15).Value = '-'
$ws.Cells.Item(8, 2).Value = '-'
$ws.Cells.Item(8, 3).Value = '-'
$ws.Cells.Item(8, 4).Value = '10:00: Expo'
$ws.Cells.Item(8, 5).Value = '-'
$ws.Cells.Item(8, 6).Value = '-'
$ws.Cells.Item(8, 7).Value = '4:00: Expo'
$ws.Cells.Item(8, 8).Value = '10:00: Expo'
$ws.Cells.Item(8, 9).Value = '4:00: Expo'
$ws.Cells.Item(8, 10).Value = '10:00: Expo'
$ws.Cells.Item(8, 11).Value = '4:00: Expo'
$ws.Cells.Item(8, 12).Value = '10:00: Expo'
$ws.Cells.Item(8, 13).Value = '4:00: Expo'
$ws.Cells.Item(8, 14).Value = '10:00: Expo'
$ws.Cells.Item(8, 15).Value = '4:00: Expo'
$ws.Cells.Item(9, 2).Value = '-'
$ws.Cells.Item(9, 3).Value = '-'
$ws.Cells.Item(9, 4).Value = '10:00: 1'
$ws.Cells.Item(9, 5).Value = '-'
$ws.Cells.Item(9, 6).Value = '-'
$ws.Cells.Item(9, 7).Value = '5:00: 1'
$ws.Cells.Item(9, 8).Value = '-'
$ws.Cells.Item(9, 9).Value = '5:00: 1'
$ws.Cells.Item(9, 10).Value = '10:00: 1'
$ws.Cells.Item(9, 11).Value = '5:00: 1'
$ws.Cells.Item(9, 12).Value = '10:00: 1'
$ws.Cells.Item(9, 13).Value = '4:00: 1'
$ws.Cells.Item(9, 14).Value = '10:00: 2'
$ws.Cells.Item(9, 15).Value = '-'
$ws.Cells.Item(10, 2).Value = '-'
$ws.Cells.Item(10, 3).Value = '-'
$ws.Cells.Item(10, 4).Value = '-'
$ws.Cells.Item(10, 5).Value = '4:00: 5'
$ws.Cells.Item(10, 6).Value = '-'
$ws.Cells.Item(10, 7).Value = '4:00: 3'
$ws.Cells.Item(10, 8).Value = '-'
$ws.Cells.Item(10, 9).Value = '-'
$ws.Cells.Item(10, 10).Value = '-'
$ws.Cells.Item(10, 11).Value = '-'
$ws.Cells.Item(10, 12).Value = '-'
$ws.Cells.Item(10, 13).Value = '-'
$ws.Cells.Item(10, 14).Value = '-'
$ws.Cells.Item(10, 15).Value = '-'
$ws.Cells.Item(11, 2).Value = '-'
$ws.Cells.Item(11, 3).Value = '-'
$ws.Cells.Item(11, 4).Value = '-'
$ws.Cells.Item(11, 5).Value = '-'
$ws.Cells.Item(11, 6).Value = '-'
$ws.Cells.Item(11, 7).Value = '-'
$ws.Cells.Item(11, 8).Value = '-'
$ws.Cells.Item(11, 9).Value = '-'
$ws.Cells.Item(11, 10).Value = '-'
$ws.Cells.Item(11, 11).Value = '-'
$ws.Cells.Item(11, 12).Value = '-'
$ws.Cells.Item(11, 13).Value = '-'
$ws.Cells.Item(11, 14).Value = '-'
$ws.Cells.Item(11, 15).Value = '-'
$ws.Cells.Item(12, 2).Value = '-'
$ws.Cells.Item(12, 3).Value = '-'
$ws.Cells.Item(12, 4).Value = '-'
$ws.Cells.Item(12, 5).Value = '-'
$ws.Cells.Item(12, 6).Value = '-'
$ws.Cells.Item(12, 7).Value = '-'
$ws.Cells.Item(12, 8).Value = '-'
$ws.Cells.Item(12, 9).Value = '4:00: 3'
$ws.Cells.Item(12, 10).Value = '-'
$ws.Cells.Item(12, 11).Value = '-'
$ws.Cells.Item(12, 12).Value = '10:00: 3'
$ws.Cells.Item(12, 13).Value = '-'
$ws.Cells.Item(12, 14).Value = '10:00: 3'
$ws.Cells.Item(12, 15).Value = '-'
$ws.Cells.Item(13, 2).Value = '-'
$ws.Cells.Item(13, 3).Value = '4:00: 5'
$ws.Cells.Item(13, 4).Value = '-'
$ws.Cells.Item(13, 5).Value = '-'
$ws.Cells.Item(13, 6).Value = '-'
$ws.Cells.Item(13, 7).Value = '-'
$ws.Cells.Item(13, 8).Value = '-'
$ws.Cells.Item(13, 9).Value = '-'
$ws.Cells.Item(13, 10).Value = '-'
$ws.Cells.Item(13, 11).Value = '-'
$ws.Cells.Item(13, 12).Value = '-'
$ws.Cells.Item(13, 13).Value = '-'
$ws.Cells.Item(13, 14).Value = '10:00: 4'
$ws.Cells.Item(13, 15).Value = '-'
$ws.Cells.Item(14, 2).Value = '-'
$ws.Cells.Item(14, 3).Value = '4:00: 1'
$ws.Cells.Item(14, 4).Value = '-'
$ws.Cells.Item(14, 5).Value = '-'
$ws.Cells.Item(14, 6).Value = '10:00: 1'
$ws.Cells.Item(14, 7).Value = '-'
$ws.Cells.Item(14, 8).Value = '-'
$ws.Cells.Item(14, 9).Value = '-'
$ws.Cells.Item(14, 10).Value = '-'
$ws.Cells.Item(14, 11).Value = '5:00: 2'
$ws.Cells.Item(14, 12).Value = '-'
$ws.Cells.Item(14, 13).Value = '-'
$ws.Cells.Item(14, 14).Value = '10:00: 1'
$ws.Cells.Item(14, 15).Value = '4:00: 1'
$ws.Cells.Item(15, 2).Value = '-'
$ws.Cells.Item(15, 3).Value = '-'
$ws.Cells.Item(15, 4).Value = '10:00: 3'
$ws.Cells.Item(15, 5).Value = '5:00: 1'
$ws.Cells.Item(15, 6).Value = '10:00: 3'
$ws.Cells.Item(15, 7).Value = '-'
$ws.Cells.Item(15, 8).Value = '10:00: 2'
$ws.Cells.Item(15, 9).Value = '-'
$ws.Cells.Item(15, 10).Value = '-'
$ws.Cells.Item(15, 11).Value = '4:00: 5'
$ws.Cells.Item(15, 12).Value = '10:00: 5'
$ws.Cells.Item(15, 13).Value = '-'
$ws.Cells.Item(15, 14).Value = '-'
$ws.Cells.Item(15, 15).Value = '-'
$ws.Cells.Item(16, 2).Value = '-'
$ws.Cells.Item(16, 3).Value = '4:00: 3'
$ws.Cells.Item(16, 4).Value = '-'
$ws.Cells.Item(16, 5).Value = '-'
$ws.Cells.Item(16, 6).Value = '-'
$ws.Cells.Item(16, 7).Value = '-'
$ws.Cells.Item(16, 8).Value = '10:00: 1'
$ws.Cells.Item(16, 9).Value = '-'
$ws.Cells.Item(16, 10).Value = '10:00: 2'
$ws.Cells.Item(16, 11).Value = '-'
$ws.Cells.Item(16, 12).Value = '-'
$ws.Cells.Item(16, 13).Value = '4:00: 5'
$ws.Cells.Item(16, 14).Value = '10:00: 5'
$ws.Cells.Item(16, 15).Value = '-'
$ws.Cells.Item(17, 2).Value = '-'
$ws.Cells.Item(17, 3).Value = '-'
$ws.Cells.Item(17, 4).Value = '-'
$ws.Cells.Item(17, 5).Value = '-'
$ws.Cells.Item(17, 6).Value = '-'
$ws.Cells.Item(17, 7).Value = '4:00: 4'
$ws.Cells.Item(17, 8).Value = '10:00: 3'
$ws.Cells.Item(17, 9).Value = '-'
$ws.Cells.Item(17, 10).Value = '-'
$ws.Cells.Item(17, 11).Value = '-'
$ws.Cells.Item(17, 12).Value = '10:00: 2'
$ws.Cells.Item(17, 13).Value = '-'
$ws.Cells.Item(17, 14).Value = '-'
$ws.Cells.Item(17, 15).Value = '-'
$ws.Cells.Item(18, 2).Value = '-'
$ws.Cells.Item(18, 3).Value = '-'
$ws.Cells.Item(18, 4).Value = '-'
$ws.Cells.Item(18, 5).Value = '-'
$ws.Cells.Item(18, 6).Value = '-'
$ws.Cells.Item(18, 7).Value = '-'
$ws.Cells.Item(18, 8).Value = '10:00: 4'
$ws.Cells.Item(18, 9).Value = '-'
$ws.Cells.Item(18, 10).Value = '-'
$ws.Cells.Item(18, 11).Value = '4:00: 4'
$ws.Cells.Item(18, 12).Value = '-'
$ws.Cells.Item(18, 13).Value = '-'
$ws.Cells.Item(18, 14).Value = '-'
$ws.Cells.Item(18, 15).Value = '-'
$ws.Cells.Item(19, 2).Value = '-'
$ws.Cells.Item(19, 3).Value = '-'
$ws.Cells.Item(19, 4).Value = '-'
$ws.Cells.Item(19, 5).Value = '-'
$ws.Cells.Item(19, 6).Value = '-'
$ws.Cells.Item(19, 7).Value = '-'
$ws.Cells.Item(19, 8).Value = '-'
$ws.Cells.Item(19, 9).Value = '-'
$ws.Cells.Item(19, 10).Value = '-'
$ws.Cells.Item(19, 11).Value = '-'
$ws.Cells.Item(19, 12).Value = '-'
$ws.Cells.Item(19, 13).Value = '-'
$ws.Cells.Item(19, 14).Value = '-'
$ws.Cells.Item(19, 15).Value = '-'
$ws.Cells.Item(20, 2).Value = '-'
$ws.Cells.Item(20, 3).Value = '4:00: 2'
$ws.Cells.Item(20, 4).Value = '10:00: 2'
$ws.Cells.Item(20, 5).Value = '-'
$ws.Cells.Item(20, 6).Value = '-'
$ws.Cells.Item(20, 7).Value = '-'
$ws.Cells.Item(20, 8).Value = '-'
$ws.Cells.Item(20, 9).Value = '-'
$ws.Cells.Item(20, 10).Value = '-'
$ws.Cells.Item(20, 11).Value = '-'
$ws.Cells.Item(20, 12).Value = '-'
$ws.Cells.Item(20, 13).Value = '-'
$ws.Cells.Item(20, 14).Value = '-'
$ws.Cells.Item(20, 15).Value = '4:00: 2'
$ws.Cells.Item(21, 2).Value = '-'
$ws.Cells.Item(21, 3).Value = '-'
$ws.Cells.Item(21, 4).Value = '10:00: 4'
$ws.Cells.Item(21, 5).Value = '4:00: 3'
$ws.Cells.Item(21, 6).Value = '10:00: 4'
$ws.Cells.Item(21, 7).Value = '-'
$ws.Cells.Item(21, 8).Value = '-'
$ws.Cells.Item(21, 9).Value = '-'
$ws.Cells.Item(21, 10).Value = '-'
$ws.Cells.Item(21, 11).Value = '-'
$ws.Cells.Item(21, 12).Value = '-'
$ws.Cells.Item(21, 13).Value = '-'
$ws.Cells.Item(21, 14).Value = '-'
$ws.Cells.Item(21, 15).Value = '4:00: 4'
$ws.Cells.Item(22, 2).Value = '-'
$ws.Cells.Item(22, 3).Value = '-'
$ws.Cells.Item(22, 4).Value = '-'
$ws.Cells.Item(22, 5).Value = '-'
$ws.Cells.Item(22, 6).Value = '-'
$ws.Cells.Item(22, 7).Value = '-'
$ws.Cells.Item(22, 8).Value = '-'
$ws.Cells.Item(22, 9).Value = '4:00: 5'
$ws.Cells.Item(22, 10).Value = '-'
$ws.Cells.Item(22, 11).Value = '-'
$ws.Cells.Item(22, 12).Value = '-'
$ws.Cells.Item(22, 13).Value = '4:00: 2'
$ws.Cells.Item(22, 14).Value = '-'
$ws.Cells.Item(22, 15).Value = '4:00: 3'
$ws.Cells.Item(23, 2).Value = '-'
$ws.Cells.Item(23, 3).Value = '-'
$ws.Cells.Item(23, 4).Value = '-'
$ws.Cells.Item(23, 5).Value = '-'
$ws.Cells.Item(23, 6).Value = '-'
$ws.Cells.Item(23, 7).Value = '-'
$ws.Cells.Item(23, 8).Value = '-'
$ws.Cells.Item(23, 9).Value = '-'
$ws.Cells.Item(23, 10).Value = '10:00: 3'
$ws.Cells.Item(23, 11).Value = '-'
$ws.Cells.Item(23, 12).Value = '-'
$ws.Cells.Item(23, 13).Value = '4:00: 4'
$ws.Cells.Item(23, 14).Value = '-'
$ws.Cells.Item(23, 15).Value = '-'
$ws.Cells.Item(24, 2).Value = '-'
$ws.Cells.Item(24, 3).Value = '-'
$ws.Cells.Item(24, 4).Value = '-'
$ws.Cells.Item(24, 5).Value = '5:00: 2'
$ws.Cells.Item(24, 6).Value = '10:00: 2'
$ws.Cells.Item(24, 7).Value = '5:00: 2'
$ws.Cells.Item(24, 8).Value = '-'
$ws.Cells.Item(24, 9).Value = '5:00: 2'
$ws.Cells.Item(24, 10).Value = '-'
$ws.Cells.Item(24, 11).Value = '4:00: 3'
$ws.Cells.Item(24, 12).Value = '10:00: 4'
$ws.Cells.Item(24, 13).Value = '4:00: 3'
$ws.Cells.Item(24, 14).Value = '-'
$ws.Cells.Item(24, 15).Value = '4:00: 5'
$ws.Cells.Item(25, 2).Value = '-'
$ws.Cells.Item(25, 3).Value = '4:00: 4'
$ws.Cells.Item(25, 4).Value = '-'
$ws.Cells.Item(25, 5).Value = '4:00: 4'
$ws.Cells.Item(25, 6).Value = '-'
$ws.Cells.Item(25, 7).Value = '4:00: 5'
$ws.Cells.Item(25, 8).Value = '-'
$ws.Cells.Item(25, 9).Value = '4:00: 4'
$ws.Cells.Item(25, 10).Value = '10:00: 4'
$ws.Cells.Item(25, 11).Value = '-'
$ws.Cells.Item(25, 12).Value = '-'
$ws.Cells.Item(25, 13).Value = '-'
$ws.Cells.Item(25, 14).Value = '-'
$ws.Cells.Item(25, 15).Value = '-'
$ws.Cells.Item(26, 2).Value = '-'
$ws.Cells.Item(26, 3).Value = '-'
$ws.Cells.Item(26, 4).Value = '11:00: H/G'
$ws.Cells.Item(26, 5).Value = '-'
$ws.Cells.Item(26, 6).Value = '-'
$ws.Cells.Item(26, 7).Value = '-'
$ws.Cells.Item(26, 8).Value = '-'
$ws.Cells.Item(26, 9).Value = '-'
$ws.Cells.Item(26, 10).Value = '-'
$ws.Cells.Item(26, 11).Value = '-'
$ws.Cells.Item(26, 12).Value = '-'
$ws.Cells.Item(26, 13).Value = '-'
$ws.Cells.Item(26, 14).Value = '11:00: H/G'
$ws.Cells.Item(26, 15).Value = '-'
$ws.Cells.Item(27, 2).Value = '-'
$ws.Cells.Item(27, 3).Value = '4:00: H/G'
$ws.Cells.Item(27, 4).Value = '-'
$ws.Cells.Item(27, 5).Value = '-'
$ws.Cells.Item(27, 6).Value = '-'
$ws.Cells.Item(27, 7).Value = '-'
$ws.Cells.Item(27, 8).Value = '-'
$ws.Cells.Item(27, 9).Value = '4:00: H/G'
$ws.Cells.Item(27, 10).Value = '11:00: H/G'
$ws.Cells.Item(27, 11).Value = '-'
$ws.Cells.Item(27, 12).Value = '-'
$ws.Cells.Item(27, 13).Value = '-'
$ws.Cells.Item(27, 14).Value = '-'
$ws.Cells.Item(27, 15).Value = '-'
$ws.Cells.Item(28, 2).Value = '-'
$ws.Cells.Item(28, 3).Value = '4:00 Runner'
$ws.Cells.Item(28, 4).Value = '-'
$ws.Cells.Item(28, 5).Value = '-'
$ws.Cells.Item(28, 6).Value = '-'
$ws.Cells.Item(28, 7).Value = '-'
$ws.Cells.Item(28, 8).Value = '-'
$ws.Cells.Item(28, 9).Value = '-'
$ws.Cells.Item(28, 10).Value = '-'
$ws.Cells.Item(28, 11).Value = '-'
$ws.Cells.Item(28, 12).Value = '11:00: H/G'
$ws.Cells.Item(28, 13).Value = '-'
$ws.Cells.Item(28, 14).Value = '-'
$ws.Cells.Item(28, 15).Value = '-'
$ws.Cells.Item(29, 2).Value = '-'
$ws.Cells.Item(29, 3).Value = '-'
$ws.Cells.Item(29, 4).Value = '-'
$ws.Cells.Item(29, 5).Value = '4:00 Runner'
$ws.Cells.Item(29, 6).Value = '-'
$ws.Cells.Item(29, 7).Value = '-'
$ws.Cells.Item(29, 8).Value = '-'
$ws.Cells.Item(29, 9).Value = '-'
$ws.Cells.Item(29, 10).Value = '-'
$ws.Cells.Item(29, 11).Value = '-'
$ws.Cells.Item(29, 12).Value = '-'
$ws.Cells.Item(29, 13).Value = '-'
$ws.Cells.Item(29, 14).Value = '-'
$ws.Cells.Item(29, 15).Value = '4:00: H/G'
$ws.Cells.Item(30, 2).Value = '-'
$ws.Cells.Item(30, 3).Value = '-'
$ws.Cells.Item(30, 4).Value = '-'
$ws.Cells.Item(30, 5).Value = '4:00: H/G'
$ws.Cells.Item(30, 6).Value = '11:00: H/G'
$ws.Cells.Item(30, 7).Value = '-'
$ws.Cells.Item(30, 8).Value = '-'
$ws.Cells.Item(30, 9).Value = '-'
$ws.Cells.Item(30, 10).Value = '-'
$ws.Cells.Item(30, 11).Value = '4:00: H/G'
$ws.Cells.Item(30, 12).Value = '-'
$ws.Cells.Item(30, 13).Value = '4:00: H/G'
$ws.Cells.Item(30, 14).Value = '-'
$ws.Cells.Item(30, 15).Value = '-'
$ws.Cells.Item(31, 2).Value = '-'
$ws.Cells.Item(31, 3).Value = '-'
$ws.Cells.Item(31, 4).Value = '-'
$ws.Cells.Item(31, 5).Value = '-'
$ws.Cells.Item(31, 6).Value = '-'
$ws.Cells.Item(31, 7).Value = '4:00: H/G'
$ws.Cells.Item(31, 8).Value = '11:00: H/G'
$ws.Cells.Item(31, 9).Value = '-'
$ws.Cells.Item(31, 10).Value = '-'
$ws.Cells.Item(31, 11).Value = '-'
$ws.Cells.Item(31, 12).Value = '-'
$ws.Cells.Item(31, 13).Value = '4:00 Runner'
$ws.Cells.Item(31, 14).Value = '-'
$ws.Cells.Item(31, 15).Value = '4:00 Runner'
$ws.Cells.Item(32, 2).Value = '-'
$ws.Cells.Item(32, 3).Value = '4:00: Expo'
$ws.Cells.Item(32, 4).Value = '-'
$ws.Cells.Item(32, 5).Value = '4:00: Expo'
$ws.Cells.Item(32, 6).Value = '10:00: Expo'
$ws.Cells.Item(32, 7).Value = '-'
$ws.Cells.Item(32, 8).Value = '-'
$ws.Cells.Item(32, 9).Value = '-'
$ws.Cells.Item(32, 10).Value = '-'
$ws.Cells.Item(32, 11).Value = '-'
$ws.Cells.Item(32, 12).Value = '-'
$ws.Cells.Item(32, 13).Value = '4:00: Expo'
$ws.Cells.Item(32, 14).Value = '4:00: Expo'
$ws.Cells.Item(32, 15).Value = '4:00: Expo'
